# Row 11 ("Fix: Game does not clear entities on 2nd entry into a map") gets
# marked Completed = Yes, with a completion date, matching the formatting
# already used by the other "Completed?" / "Completion Date" cells on the
# sheet (e.g. row 6 / row 12: highlighted fill for C, m/d/yyyy date format
# for D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the existing formatting from a sibling "done" row (row 6) onto the
# new cells so the highlighted fill (Completed? column) and date number
# format (Completion Date column) match exactly.
$ws.Range("C6").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D11").PasteSpecial(-4122)

# Fill in the values.
$ws.Range("C11").Value = "Yes"
$ws.Range("D11").Value = 45472

# Reselect the whole used range, mirroring the saved selection state.
$ws.Range("A1:D32").Select()
